# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 7656
$wsExhibit.Range("F6").Value = 5567
$wsExhibit.Range("F7").Value = 456
$wsExhibit.Range("F10").Value = 66
$wsExhibit.Range("F12").Value = 195
$wsExhibit.Range("F13").Value = 49

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 7656
$wsAll.Range("F6").Value = 5567
$wsAll.Range("F7").Value = 456
$wsAll.Range("F10").Value = 66
$wsAll.Range("F14").Value = 195
$wsAll.Range("F15").Value = 49
